$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the computed ellipse() sample strings in column I to reflect the
# current x/y/diam values for each row.
$ws.Range("I2").Value = "ellipse(562,419,16,16)"
$ws.Range("I3").Value = "ellipse(162,389,29,29)"

# Move the active selection to D11.
[void]$ws.Range("D11").Select()
